$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Remove the "Square" formula cell at G3 (C3*C3 column is being dropped from row 3)
[void]$ws.Range("G3").ClearContents()

# New row 5: a "Ken" entry added below the existing Joe / Mary rows
$ws.Range("A5").Value = 42037
$ws.Range("A5").NumberFormat = "DD/MM/YY"
$ws.Range("B5").Value = "Ken"
$ws.Range("C5").Value = 900
$ws.Range("D5").Value = 100
$ws.Range("E5").Formula = "=D5/C5"
$ws.Range("G5").Value = 30

$ws.Rows("5").RowHeight = 13.55

# Move the selection like the author left it
[void]$ws.Range("G4").Select()
